# Apply the "pricing_init_rate_period" / "pricing_prepenalty_allowed" /
# "pricing_prepenalty_exists" column insertion to the "invalid" sheet.
#
# Strategy: column R (pricing_fixed_rate) and the old column S
# (pricing_var_margin) need to shift right by one column to make room for
# a brand-new column R (pricing_init_rate_period). Two more brand-new
# columns are appended at U and V. To keep the shared-string table (and
# the final cell contents) correct we first read every old value that is
# going to move, and only then write all the new values, left to right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "invalid" sheet

function Set-PlainValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-EmptyCell($row, $col) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "TEMP"
    $cell.ClearContents()
    $cell.Style = "Normal"
}

# Columns (1-based): Q=17 R=18 S=19 T=20 U=21 V=22

# ---- Row 1 (headers) ----
$oldR1 = $ws.Cells.Item(1, 18).Value()
$oldS1 = $ws.Cells.Item(1, 19).Value()
Set-PlainValue 1 18 "pricing_init_rate_period"
$ws.Cells.Item(1, 19).Value = $oldR1
$ws.Cells.Item(1, 20).Value = $oldS1
Set-PlainValue 1 21 "pricing_prepenalty_allowed"
Set-PlainValue 1 22 "pricing_prepenalty_exists"

# ---- Row 2 ----
$oldS2 = $ws.Cells.Item(2, 19).Value()
Set-PlainValue 2 18 24
$ws.Cells.Item(2, 20).Value = $oldS2
Set-PlainValue 2 21 1
Set-PlainValue 2 22 1

# ---- Row 3 ----
$oldR3 = $ws.Cells.Item(3, 18).Value()
$oldS3 = $ws.Cells.Item(3, 19).Value()
Set-PlainValue 3 18 36
$ws.Cells.Item(3, 19).Value = $oldR3
$ws.Cells.Item(3, 20).Value = $oldS3
Set-PlainValue 3 21 2
Set-PlainValue 3 22 2

# ---- Row 4 ----
$oldR4 = $ws.Cells.Item(4, 18).Value()
$oldS4 = $ws.Cells.Item(4, 19).Value()
Set-PlainValue 4 18 1
$ws.Cells.Item(4, 19).Value = $oldR4
$ws.Cells.Item(4, 20).Value = $oldS4
Set-PlainValue 4 21 999
Set-PlainValue 4 22 1

# ---- Row 5 ----
$oldR5 = $ws.Cells.Item(5, 18).Value()
$oldS5 = $ws.Cells.Item(5, 19).Value()
Set-PlainValue 5 18 2
$ws.Cells.Item(5, 19).Value = $oldR5
$ws.Cells.Item(5, 20).Value = $oldS5
Set-PlainValue 5 21 0
Set-PlainValue 5 22 2

# ---- Row 6 ----
$oldR6 = $ws.Cells.Item(6, 18).Value()
$oldS6 = $ws.Cells.Item(6, 19).Value()
Set-PlainValue 6 18 2.5
$ws.Cells.Item(6, 19).Value = $oldR6
$ws.Cells.Item(6, 20).Value = $oldS6
Set-PlainValue 6 21 3
Set-PlainValue 6 22 0

# ---- Row 7 ----
$oldR7 = $ws.Cells.Item(7, 18).Value()
Set-PlainValue 7 18 0
$ws.Cells.Item(7, 19).Value = $oldR7
Set-PlainValue 7 21 1
Set-PlainValue 7 22 0

# ---- Row 8 ----
Set-PlainValue 8 18 1
Set-EmptyCell 8 20
Set-PlainValue 8 21 2
Set-PlainValue 8 22 999

# ---- Row 9 ----
Set-PlainValue 9 18 5
Set-PlainValue 9 21 999
Set-PlainValue 9 22 999

# ---- Row 10 ----
Set-EmptyCell 10 18
Set-PlainValue 10 21 0
Set-PlainValue 10 22 1

# ---- Row 11 ----
Set-PlainValue 11 18 7
Set-PlainValue 11 21 0
Set-PlainValue 11 22 2

# ---- Column widths ----
$ws.Columns.Item(18).ColumnWidth = 20.1666666   # -> width 21
$ws.Columns.Item(21).ColumnWidth = 23.25        # -> width ~24.1667 (closest achievable to 24.1640625)
$ws.Columns.Item(22).ColumnWidth = 21.5         # -> width ~22.3333 (closest achievable to 22.33203125)

# ---- View / selection ----
$ws.Range("T1").Select()
